$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newValues = @(
    "34-20=14",
    "46-29=17",
    "24+64=88",
    "29+27=56",
    "24+14=38",
    "17+63=80",
    "83-29=54",
    "74+17=91",
    "45+24=69",
    "75-22=53",
    "30+57=87",
    "66-53=13",
    "71-41=30",
    "28+58=86",
    "44+15=59",
    "63+7=70",
    "57+15=72",
    "86+1=87",
    "77-17=60",
    "41+55=96",
    "45-11=34",
    "23-20=3",
    "47+44=91",
    "0+64=64",
    "64+2=66",
    "57+20=77",
    "83-55=28",
    "98-57=41",
    "40-3=37",
    "95-89=6",
    "78-8=70",
    "16-11=5",
    "83-0=83",
    "68-56=12",
    "29+42=71",
    "71+21=92",
    "73-3=70",
    "39-35=4",
    "74+9=83",
    "34-0=34",
    "2+23=25",
    "81-59=22",
    "22+65=87",
    "21+32=53",
    "83-58=25",
    "74-31=43",
    "56+29=85",
    "45-4=41",
    "52-51=1",
    "0+87=87",
    "84-1=83",
    "47+34=81",
    "58-19=39",
    "28-8=20",
    "17+46=63",
    "97-1=96",
    "67-18=49",
    "4+78=82",
    "63+6=69",
    "16+13=29",
    "82-33=49",
    "67+9=76",
    "25+47=72",
    "69-34=35",
    "49+29=78",
    "27+47=74",
    "71-25=46",
    "69-15=54",
    "60-31=29",
    "17+77=94",
    "17+35=52",
    "81-76=5",
    "19+6=25",
    "97-11=86",
    "79-47=32",
    "26+28=54",
    "95-58=37",
    "91-2=89",
    "48+3=51",
    "61-16=45",
    "90+1=91",
    "14+43=57",
    "29-16=13",
    "79-5=74",
    "25+70=95",
    "60-23=37",
    "61+27=88",
    "0+41=41",
    "90-43=47",
    "33-2=31",
    "93-59=34",
    "35-9=26",
    "96-25=71",
    "44+3=47",
    "70+9=79",
    "12+27=39",
    "5+85=90",
    "23+71=94",
    "78-14=64",
    "3+36=39"
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."
